$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Text
    )
    $cell = $ws.Range($Address)
    # Force text interpretation so values like "131.35" or "8.00" are not
    # silently coerced into numbers by COM's Value setter, then restore the
    # cell's original (default) style so no stray per-cell number format
    # sticks around.
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "60.589.93"
Set-TextValue "E2" "  +0.47%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.337.95"
Set-TextValue "E3" "  -0.12%  "

# Row 4 - TetherUSD
Set-TextValue "E4" "  +0.01%  "

# Row 5 - BNB
Set-TextValue "D5" "552.13"
Set-TextValue "E5" "  +1.14%  "

# Row 6 - Solana
Set-TextValue "D6" "131.35"
Set-TextValue "E6" "  -0.39%  "

# Row 8 - XRP
Set-TextValue "E8" "  -0.71%  "

# Row 9 - LidoStakedEther
Set-TextValue "D9" "2.336.14"
Set-TextValue "E9" "  -0.06%  "

# Row 10 - Dogecoin
Set-TextValue "E10" "  +1.15%  "

# Row 11 - Toncoin
Set-TextValue "E11" "  +1.68%  "

# Row 12 - TRON
Set-TextValue "E12" "  -0.45%  "

# Row 13 - Cardano
Set-TextValue "E13" "  +1.38%  "

# Row 14 - Avalanche
Set-TextValue "D14" "23.96"
Set-TextValue "E14" "  +0.54%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "2.755.40"
Set-TextValue "E15" "  -0.03%  "

# Row 16 - WrappedBTC
Set-TextValue "D16" "60.499.18"
Set-TextValue "E16" "  +0.41%  "

# Row 17 - ShibaInu
Set-TextValue "E17" "  +1.11%  "

# Row 18 - WrappedEther
Set-TextValue "D18" "2.318.58"
Set-TextValue "E18" "  -1.11%  "

# Row 20 - Polkadot
Set-TextValue "D20" "4.11"
Set-TextValue "E20" "  -0.99%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "314.97"
Set-TextValue "E21" "  +0.37%  "

# Row 22 - Uniswap
Set-TextValue "D22" "6.63"
Set-TextValue "E22" "  -2.48%  "

# Row 24 - Litecoin
Set-TextValue "D24" "64.24"
Set-TextValue "E24" "  +1.17%  "

# Row 25 - Kaspa
Set-TextValue "E25" "  -0.49%  "

# Row 26 - Binance-PegBSC-USD
Set-TextValue "E26" "  -0.16%  "

# Row 27 - InternetComputer(DFINITY)
Set-TextValue "D27" "8.00"
Set-TextValue "E27" "  +1.07%  "

# Row 28 - Fetch.AI
Set-TextValue "E28" "  +3.28%  "

# Row 29 - SuiNetwork
Set-TextValue "E29" "  +9.20%  "

# Row 30 - PancakeSwap
Set-TextValue "E30" "  +0.07%  "

# Row 31 - Monero
Set-TextValue "D31" "171.30"
Set-TextValue "E31" "  -0.23%  "

# Row 32 - PEPE
Set-TextValue "E32" "  +0.68%  "

# Row 33 - Aptos
Set-TextValue "E33" "  +2.51%  "

# Row 34 - PolygonEcosystemToken
Set-TextValue "D34" "0.386"
Set-TextValue "E34" "  +0.81%  "

# Row 35 - ImmutableX
Set-TextValue "E35" "  -1.13%  "

# Row 36 - EthereumClassic
Set-TextValue "D36" "18.08"
Set-TextValue "E36" "  +0.38%  "

# Row 38 - FirstDigitalUSD
Set-TextValue "E38" "  -0.02%  "

# Row 39 - NEARProtocol
Set-TextValue "E39" "  -0.19%  "

# Row 40 - Bittensor
Set-TextValue "D40" "331.83"
Set-TextValue "E40" "  +3.19%  "

# Row 41 - Stacks
Set-TextValue "E41" "  +0.16%  "

# Row 42 - OKB
Set-TextValue "D42" "38.10"
Set-TextValue "E42" "  -0.12%  "

# Row 43 - Aave
Set-TextValue "D43" "138.74"
Set-TextValue "E43" "  -1.68%  "

# Row 44 - Filecoin
Set-TextValue "E44" "  +1.82%  "

# Row 45 - Stellar
Set-TextValue "D45" "0.0951"
Set-TextValue "E45" "  +0.63%  "

# Row 46 - InjectiveProtocol
Set-TextValue "D46" "19.36"
Set-TextValue "E46" "  -0.78%  "

# Row 47 - Mantle
Set-TextValue "D47" "0.569"
Set-TextValue "E47" "  +1.63%  "

# Row 48 - Hedera
Set-TextValue "D48" "0.0500"
Set-TextValue "E48" "  +0.79%  "

# Row 49 - BabyDogeCoin
Set-TextValue "D49" "0.0₆0226"
Set-TextValue "E49" "  +7.71%  "

# Row 50 - VeChain
Set-TextValue "E50" "  +1.69%  "

# Row 51 - WhiteBITCoin -> EnergySwap
Set-TextValue "B51" "EnergySwap"
Set-TextValue "C51" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D51" "17.18"
Set-TextValue "E51" "  +1.73%  "
